# Khulo disability_prevalence.xlsx update
#
# Summary of the edit (per the source-control diff):
#  - New report title in row 1 (merged A1:I1), taller row, bold centred Arial 11.
#  - Row 4 used to be a label + a row of "..." placeholders; it now carries the
#    "family with disabilities" series with real counts (1093 ... 1056).
#  - Row 5 used to hold the "Source:" note (merged A5:H5); it now carries the
#    "disabilities" series with real counts (1339 ... 1309).
#  - The "Source:" note moves down to row 6 (merged A6:H6), replacing the old
#    "Note: confidential/unavailable" remark, which is dropped entirely.
#  - Column A narrows a bit, row heights/number formats are touched up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Break the old merge (old "Source" row) before values/rows are reshuffled
# ---------------------------------------------------------------------------
$ws.Range("A5:H5").UnMerge()

# ---------------------------------------------------------------------------
# 1. Row 1 - new title, merged across A1:I1
# ---------------------------------------------------------------------------
$titleRange = $ws.Range("A1:I1")
$titleRange.Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Khulo Municipality"
$titleRange.Font.Name = "Arial"
$titleRange.Font.Size = 11
$titleRange.Font.Bold = $true
$titleRange.Font.Underline = $false
$titleRange.HorizontalAlignment = "xlCenter"
$titleRange.VerticalAlignment = "xlCenter"
$titleRange.WrapText = $true
$titleRange.Borders.Item("xlEdgeTop").LineStyle = "xlLineStyleNone"
$titleRange.Borders.Item("xlEdgeBottom").LineStyle = "xlLineStyleNone"
$ws.Rows.Item(1).RowHeight = 51
$titleRange.Merge()

# ---------------------------------------------------------------------------
# 2. Row 2 - unchanged text, but no longer a custom row height
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# 3. Row 3 - the (blank) A3 cell now inherits the Sylfaen column font
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 4. Row 4 -> "family with disabilities Persons" series
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "

$dataCols = @("B","C","D","E","F","G","H","I")
$row4Values = @(1093,1099,1071,1104,1113,1108,1053,1056)
$row4Range = $ws.Range("B4:I4")
for ($i = 0; $i -lt $dataCols.Length; $i++) {
    $ws.Range($dataCols[$i] + "4").Value = $row4Values[$i]
}
$row4Range.NumberFormat = "#\ ##0"
$row4Range.Font.Name = "Arial"
$row4Range.Font.Size = 10
$row4Range.HorizontalAlignment = "xlGeneral"
$row4Range.VerticalAlignment = "xlBottom"
$row4Range.WrapText = $false
$row4Range.Borders.Item("xlEdgeTop").LineStyle = "xlLineStyleNone"
$row4Range.Borders.Item("xlEdgeBottom").LineStyle = "xlLineStyleNone"
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 5. Row 5 -> "disabilities Persons" series
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").Font.Underline = $false
$ws.Range("A5").Borders.Item("xlEdgeTop").LineStyle = "xlLineStyleNone"
$ws.Range("A5").Borders.Item("xlEdgeBottom").LineStyle = "xlContinuous"

$row5Values = @(1339,1345,1330,1377,1379,1376,1312,1309)
for ($i = 0; $i -lt $dataCols.Length; $i++) {
    $ws.Range($dataCols[$i] + "5").Value = $row5Values[$i]
}
$row5Range = $ws.Range("B5:I5")
$row5Range.NumberFormat = "#\ ##0"
$row5Range.Font.Name = "Arial"
$row5Range.Font.Size = 10
$row5Range.HorizontalAlignment = "xlGeneral"
$row5Range.VerticalAlignment = "xlBottom"
$row5Range.WrapText = $false
$row5Range.Borders.Item("xlEdgeTop").LineStyle = "xlLineStyleNone"
$row5Range.Borders.Item("xlEdgeBottom").LineStyle = "xlLineStyleNone"
# last cell of the row keeps a bottom rule (matches the header rule above it)
$ws.Range("I5").Borders.Item("xlEdgeBottom").LineStyle = "xlContinuous"
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 6. Row 6 -> "Source: ..." note (was row 5's merged note; "Note:" row removed)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$noteRange = $ws.Range("A6:H6")
$noteRange.Font.Name = "Arial"
$noteRange.Font.Size = 9
$noteRange.Font.Bold = $true
$noteRange.Font.Underline = $true
$noteRange.HorizontalAlignment = "xlLeft"
$noteRange.VerticalAlignment = "xlCenter"
$noteRange.WrapText = $true
$ws.Range("B6:H6").Borders.Item("xlEdgeTop").LineStyle = "xlContinuous"
$ws.Range("A6").Borders.Item("xlEdgeTop").LineStyle = "xlLineStyleNone"
$ws.Rows.Item(6).RowHeight = 27.75
$noteRange.Merge()

# ---------------------------------------------------------------------------
# 7. Column A width
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20

# ---------------------------------------------------------------------------
# 8. Selection, matching the saved view state
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Select()

Write-Host "edit applied"
